# Trade #50 closed at 2026-02-16 21:30:30 - leadlag UP +0.000%
#
# This edit:
#   1. Closes leadlag Trade #20 (row 19) -> CLOSED with exit price/P&L/duration.
#   2. Mirrors that closed trade into the "All Trades" sheet as a new row.
#   3. Appends a brand-new OPEN leadlag Trade #50 (row 40).
#   4. Refreshes the Summary and Comparison roll-up statistics.
#
# NOTE: several "numeric looking" strings (dates, times, percentages, and
# decimal-looking ratios) are stored in the source workbook as literal TEXT,
# not numbers. Assigning a plain string to Range.Value would otherwise be
# auto-coerced to a Number/Date/Percentage by this engine (mirroring real
# Excel's smart entry), so for those cells we force the NumberFormat to
# Text ("@") immediately before the assignment.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("C2").Value = 20

$summary.Range("D2").NumberFormat = "@"
$summary.Range("D2").Value = "65.0%"

$summary.Range("E2").NumberFormat = "@"
$summary.Range("E2").Value = "+3.5966%"

$summary.Range("F2").NumberFormat = "@"
$summary.Range("F2").Value = "+0.1798%"

$summary.Range("C3").Value = 38

$summary.Range("D3").NumberFormat = "@"
$summary.Range("D3").Value = "31.6%"

$summary.Range("E3").NumberFormat = "@"
$summary.Range("E3").Value = "+3.5513%"

$summary.Range("F3").NumberFormat = "@"
$summary.Range("F3").Value = "+0.0935%"

# ---------------------------------------------------------------------
# leadlag sheet - close Trade #20 (row 19)
# ---------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

$leadlag.Range("G19").Value = 68412.082132

$leadlag.Range("H19").Value = "CLOSED"

$leadlag.Range("I19").Value = 1.0574
$leadlag.Range("J19").Value = 10.57

$leadlag.Range("M19").Value = "time_exit_5min"

$leadlag.Range("N19").Value = 5

# ---------------------------------------------------------------------
# leadlag sheet - append new OPEN Trade #50 (row 40)
# ---------------------------------------------------------------------
$leadlag.Range("A40").Value = 50

$leadlag.Range("B40").NumberFormat = "@"
$leadlag.Range("B40").Value = "2026-02-16"

$leadlag.Range("C40").NumberFormat = "@"
$leadlag.Range("C40").Value = "21:30:30"

$leadlag.Range("D40").Value = "leadlag"
$leadlag.Range("E40").Value = "UP"

$leadlag.Range("F40").Value = 68818.08

# G40 intentionally left blank - trade is still OPEN, no exit price yet.

$leadlag.Range("H40").Value = "OPEN"

$leadlag.Range("I40").Value = 0
$leadlag.Range("J40").Value = 0
$leadlag.Range("K40").Value = 0.75

$leadlag.Range("L40").Value = "Binance leading with 0.153% move"

# M40 intentionally left blank - no exit reason yet.

$leadlag.Range("N40").Value = 0

# ---------------------------------------------------------------------
# All Trades sheet - mirror the newly CLOSED Trade #20 as row 21
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("A21").Value = 20

$allTrades.Range("B21").NumberFormat = "@"
$allTrades.Range("B21").Value = "2026-02-16"

$allTrades.Range("C21").NumberFormat = "@"
$allTrades.Range("C21").Value = "21:25:29"

$allTrades.Range("D21").Value = "leadlag"
$allTrades.Range("E21").Value = "DOWN"

$allTrades.Range("F21").Value = 69143.23
$allTrades.Range("G21").Value = 68412.082132

$allTrades.Range("H21").Value = "CLOSED"

$allTrades.Range("I21").Value = 1.0574
$allTrades.Range("J21").Value = 10.57
$allTrades.Range("K21").Value = 0.75

$allTrades.Range("L21").Value = "Binance leading with -0.155% move"

$allTrades.Range("M21").Value = "time_exit_5min"

$allTrades.Range("N21").Value = 5

# ---------------------------------------------------------------------
# Comparison sheet
# ---------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Range("B2").Value = 38

$comparison.Range("C2").NumberFormat = "@"
$comparison.Range("C2").Value = "31.6%"

$comparison.Range("D2").NumberFormat = "@"
$comparison.Range("D2").Value = "2.54"

$comparison.Range("E2").NumberFormat = "@"
$comparison.Range("E2").Value = "+0.4887%"

$comparison.Range("G2").NumberFormat = "@"
$comparison.Range("G2").Value = "1.27"
